$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "xnVmb914"
$ws.Range("B2").Value = 23090109
$ws.Range("C2").Value = "mwuqrpx90"
$ws.Range("D2").Value = "BxY7e2#&"
$ws.Range("F2").Value = "VTlXAAGy"
$ws.Range("G2").Value = "KMij"
